$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1630846444411361
$ws.Cells.Item(2, 3).Value = -1.088872551767469
$ws.Cells.Item(2, 4).Value = 0.2388572144868575
$ws.Cells.Item(2, 5).Value = -0.1694399950501343
$ws.Cells.Item(2, 6).Value = 0.5040425926957953
$ws.Cells.Item(2, 7).Value = -0.1932960119371557
$ws.Cells.Item(2, 8).Value = 0.04367710249015017
$ws.Cells.Item(2, 9).Value = 0.07925193655196383
$ws.Cells.Item(2, 10).Value = 0.9052087946043501
$ws.Cells.Item(2, 11).Value = 0.2491603105793238
$ws.Cells.Item(3, 2).Value = 0.6966127121850244
$ws.Cells.Item(3, 3).Value = 0.03971143396413679
$ws.Cells.Item(3, 4).Value = 0.6050679262073962
$ws.Cells.Item(3, 5).Value = -0.1390344077924999
$ws.Cells.Item(3, 6).Value = 0.07771993578465181
$ws.Cells.Item(3, 7).Value = 0.1045584651549948
$ws.Cells.Item(3, 8).Value = 0.9267429767877318
$ws.Cells.Item(3, 9).Value = 0.2690667718748542
$ws.Cells.Item(3, 10).Value = -0.3899018488887417
$ws.Cells.Item(3, 11).Value = 0.008856469899756692
$ws.Cells.Item(4, 2).Value = 0.5889366529896609
$ws.Cells.Item(4, 3).Value = -0.1629753181314043
$ws.Cells.Item(4, 4).Value = 0.05339438235456767
$ws.Cells.Item(4, 5).Value = 0.08260390735964929
$ws.Cells.Item(4, 6).Value = 0.906862924893769
$ws.Cells.Item(4, 7).Value = 0.2504397716903112
$ws.Cells.Item(4, 8).Value = -0.4078772338215886
$ws.Cells.Item(4, 9).Value = -0.00880507520269791
$ws.Cells.Item(4, 10).Value = -0.0702522111134038
$ws.Cells.Item(4, 11).Value = -0.5044206789217187
$ws.Cells.Item(5, 2).Value = 0.135347353279659
$ws.Cells.Item(5, 3).Value = 0.0918012475915978
$ws.Cells.Item(5, 4).Value = 0.8871553623927978
$ws.Cells.Item(5, 5).Value = 0.2212014561932883
$ws.Cells.Item(5, 6).Value = -0.4400216186372514
$ws.Cells.Item(5, 7).Value = -0.04180104165272297
$ws.Cells.Item(5, 8).Value = -0.1034922560287353
$ws.Cells.Item(5, 9).Value = -0.5377297743012852
$ws.Cells.Item(5, 10).Value = 0.4410516555687047
$ws.Cells.Item(5, 11).Value = 0.2139302114848279
$ws.Cells.Item(6, 2).Value = 0.8583540716299061
$ws.Cells.Item(6, 3).Value = 0.219077672514715
$ws.Cells.Item(6, 4).Value = -0.4327386538548288
$ws.Cells.Item(6, 5).Value = -0.03138452870947878
$ws.Cells.Item(6, 6).Value = -0.09205111059441351
$ws.Cells.Item(6, 7).Value = -0.5259549742272362
$ws.Cells.Item(6, 8).Value = 0.452935276392097
$ws.Cells.Item(6, 9).Value = 0.2258494563263944
$ws.Cells.Item(6, 10).Value = -0.4076000933034996
$ws.Cells.Item(6, 11).Value = 0.0198003091970762
$ws.Cells.Item(7, 2).Value = -0.5546208149506421
$ws.Cells.Item(7, 3).Value = -0.1232296976428027
$ws.Cells.Item(7, 4).Value = -0.1719238571658218
$ws.Cells.Item(7, 5).Value = -0.6008843017759042
$ws.Cells.Item(7, 6).Value = 0.3800526390729398
$ws.Cells.Item(7, 7).Value = 0.1538132908583271
$ws.Cells.Item(7, 8).Value = -0.4792866254931878
$ws.Cells.Item(7, 9).Value = -0.05174196790978589
$ws.Cells.Item(7, 10).Value = -0.2180182939146202
$ws.Cells.Item(7, 11).Value = -0.09528239010531198
$ws.Cells.Item(8, 2).Value = -0.07660908550845907
$ws.Cells.Item(8, 3).Value = -0.5226165970886194
$ws.Cells.Item(8, 4).Value = 0.4506431296588137
$ws.Cells.Item(8, 5).Value = 0.2208746131534756
$ws.Cells.Item(8, 6).Value = -0.4138532542967197
$ws.Cells.Item(8, 7).Value = 0.01293740259710929
$ws.Cells.Item(8, 8).Value = -0.153689426264273
$ws.Cells.Item(8, 9).Value = -0.03111699931576278
$ws.Cells.Item(8, 10).Value = 0.01605047729302533
$ws.Cells.Item(8, 11).Value = -0.3805999332580445
$ws.Cells.Item(9, 2).Value = 0.7107656714615536
$ws.Cells.Item(9, 3).Value = 0.3617805447359003
$ws.Cells.Item(9, 4).Value = -0.3284087146615119
$ws.Cells.Item(9, 5).Value = 0.07258913989446059
$ws.Cells.Item(9, 6).Value = -0.106032695953167
$ws.Cells.Item(9, 7).Value = 0.01096152432303771
$ws.Cells.Item(9, 8).Value = 0.0555349373482778
$ws.Cells.Item(9, 9).Value = -0.3423217815153134
$ws.Cells.Item(9, 10).Value = -0.0888458464935532
$ws.Cells.Item(9, 11).Value = 0.1310104539444903
$ws.Cells.Item(10, 2).Value = -0.5260169647433546
$ws.Cells.Item(10, 3).Value = -0.05137068376835741
$ws.Cells.Item(10, 4).Value = -0.1974123165853111
$ws.Cells.Item(10, 5).Value = -0.06599669411565129
$ws.Cells.Item(10, 6).Value = -0.01502993377330319
$ws.Cells.Item(10, 7).Value = -0.410047744940668
$ws.Cells.Item(10, 8).Value = -0.1553091497739889
$ws.Cells.Item(10, 9).Value = 0.06510967366033998
$ws.Cells.Item(10, 10).Value = -0.1705525061135018
$ws.Cells.Item(10, 11).Value = -0.3494711127835084
$ws.Cells.Item(11, 2).Value = -0.1300102632924844
$ws.Cells.Item(11, 3).Value = 0.00888099880676374
$ws.Cells.Item(11, 4).Value = 0.06313616294635133
$ws.Cells.Item(11, 5).Value = -0.3304359561846383
$ws.Cells.Item(11, 6).Value = -0.07506199009359621
$ws.Cells.Item(11, 7).Value = 0.1456359881641034
$ws.Cells.Item(11, 8).Value = -0.08990357932170889
$ws.Cells.Item(11, 9).Value = -0.2687683468051136
$ws.Cells.Item(11, 10).Value = -0.2792820200016274
$ws.Cells.Item(11, 11).Value = -0.4869807987902649
$ws.Cells.Item(12, 2).Value = 0.05077946203237782
$ws.Cells.Item(12, 3).Value = -0.344626946110449
$ws.Cells.Item(12, 4).Value = -0.09000824431623794
$ws.Cells.Item(12, 5).Value = 0.1303795298704019
$ws.Cells.Item(12, 6).Value = -0.1052864881356162
$ws.Cells.Item(12, 7).Value = -0.2842023355090271
$ws.Cells.Item(12, 8).Value = -0.2947364119262331
$ws.Cells.Item(12, 9).Value = -0.5024432253954677
$ws.Cells.Item(12, 10).Value = -0.2852286099811868
$ws.Cells.Item(12, 11).Value = 0.02446021624927014
$ws.Cells.Item(13, 2).Value = 0.0822944582731665
$ws.Cells.Item(13, 3).Value = 0.2248716626294139
$ws.Cells.Item(13, 4).Value = -0.04696104924847117
$ws.Cells.Item(13, 5).Value = -0.2426806063420137
$ws.Cells.Item(13, 6).Value = -0.2610220372837692
$ws.Cells.Item(13, 7).Value = -0.4723563104746108
$ws.Cells.Item(13, 8).Value = -0.2568270884946112
$ws.Cells.Item(13, 9).Value = 0.05207866879825723
$ws.Cells.Item(13, 10).Value = -0.1025880601374284
$ws.Cells.Item(13, 11).Value = 0.5476698920998867
$ws.Cells.Item(14, 2).Value = -0.1414695961309799
$ws.Cells.Item(14, 3).Value = -0.2791950053566942
$ws.Cells.Item(14, 4).Value = -0.2708243371007503
$ws.Cells.Item(14, 5).Value = -0.4698642071905083
$ws.Cells.Item(14, 6).Value = -0.2486761060184952
$ws.Cells.Item(14, 7).Value = 0.06283442082819635
$ws.Cells.Item(14, 8).Value = -0.09063329381595758
$ws.Cells.Item(14, 9).Value = 0.560176603242597
$ws.Cells.Item(14, 10).Value = 0.3441873680460773
$ws.Cells.Item(14, 11).Value = -0.06276115523468917
$ws.Cells.Item(15, 2).Value = -0.1284196148407921
$ws.Cells.Item(15, 3).Value = -0.3892880536239818
$ws.Cells.Item(15, 4).Value = -0.1956617690156243
$ws.Cells.Item(15, 5).Value = 0.1035623046680138
$ws.Cells.Item(15, 6).Value = -0.05538318383206281
$ws.Cells.Item(15, 7).Value = 0.5929841687531761
$ws.Cells.Item(15, 8).Value = 0.3759056446645694
$ws.Cells.Item(15, 9).Value = -0.03152873430779118
$ws.Cells.Item(15, 10).Value = 0.4673336997040285
$ws.Cells.Item(15, 11).Value = 0.3084050283952137
$ws.Cells.Item(16, 2).Value = -0.01772937608965819
$ws.Cells.Item(16, 3).Value = 0.1904199002428446
$ws.Cells.Item(16, 4).Value = -0.01245745804686382
$ws.Cells.Item(16, 5).Value = 0.6148928434241717
$ws.Cells.Item(16, 6).Value = 0.3877598625592121
$ws.Cells.Item(16, 7).Value = -0.02448210058422801
$ws.Cells.Item(16, 8).Value = 0.4720825582591874
$ws.Cells.Item(16, 9).Value = 0.3120561101207888
$ws.Cells.Item(16, 10).Value = 0.3676602717569419
$ws.Cells.Item(16, 11).Value = 2.617197116239002
$ws.Cells.Item(17, 2).Value = 0.2009879203755026
$ws.Cells.Item(17, 3).Value = -0.003682669475092215
$ws.Cells.Item(17, 4).Value = 0.6226952474679071
$ws.Cells.Item(17, 5).Value = 0.3950699233937206
$ws.Cells.Item(17, 6).Value = -0.01742125293059998
$ws.Cells.Item(17, 7).Value = 0.4790177347602291
$ws.Cells.Item(17, 8).Value = 0.3189281133539701
$ws.Cells.Item(17, 9).Value = 0.3745006084643348
$ws.Cells.Item(17, 10).Value = 2.624021620117335
$ws.Cells.Item(17, 11).Value = 10.09303763079093
$ws.Cells.Item(18, 2).Value = -0.1106570732907728
$ws.Cells.Item(18, 3).Value = 0.5574356245296381
$ws.Cells.Item(18, 4).Value = 0.3494793248607413
$ws.Cells.Item(18, 5).Value = -0.0537519294468537
$ws.Cells.Item(18, 6).Value = 0.4470476400263193
$ws.Cells.Item(18, 7).Value = 0.2890118858004561
$ws.Cells.Item(18, 8).Value = 0.3455519702652192
$ws.Cells.Item(18, 9).Value = 2.59552891386994
$ws.Cells.Item(18, 10).Value = 10.06475980582587
$ws.Cells.Item(18, 11).Value = -8.106464717219385
$ws.Cells.Item(19, 2).Value = 0.5913088222741343
$ws.Cells.Item(19, 3).Value = 0.3522324851175385
$ws.Cells.Item(19, 4).Value = -0.06534086029760827
$ws.Cells.Item(19, 5).Value = 0.4289262259694052
$ws.Cells.Item(19, 6).Value = 0.267911772437701
$ws.Cells.Item(19, 7).Value = 0.3230936335787457
$ws.Cells.Item(19, 8).Value = 2.57245120070562
$ws.Cells.Item(19, 9).Value = 10.04139962167697
$ws.Cells.Item(19, 10).Value = -8.129953734994666
$ws.Cells.Item(19, 11).Value = 0.009676295449283501
$ws.Cells.Item(20, 2).Value = 0.09495705561128835
$ws.Cells.Item(20, 3).Value = -0.1966201317798049
$ws.Cells.Item(20, 4).Value = 0.3555705446956575
$ws.Cells.Item(20, 5).Value = 0.2206704977705274
$ws.Cells.Item(20, 6).Value = 0.2876467492000076
$ws.Cells.Item(20, 7).Value = 2.542330314111783
$ws.Cells.Item(20, 8).Value = 10.01368383765935
$ws.Cells.Item(20, 9).Value = -8.156583427754628
$ws.Cells.Item(20, 10).Value = -0.01646294162652273
$ws.Cells.Item(20, 11).Value = 2.131310184089503
$ws.Cells.Item(21, 2).Value = -0.2398842457485278
$ws.Cells.Item(21, 3).Value = 0.3272291486687826
$ws.Cells.Item(21, 4).Value = 0.1959944722584051
$ws.Cells.Item(21, 5).Value = 0.2633601162231758
$ws.Cells.Item(21, 6).Value = 2.517707017856339
$ws.Cells.Item(21, 7).Value = 9.988700664718653
$ws.Cells.Item(21, 8).Value = -8.181813387285533
$ws.Cells.Item(21, 9).Value = -0.04183851343822592
$ws.Cells.Item(21, 10).Value = 2.105855004056436
$ws.Cells.Item(21, 11).Value = -1.334125068513618
$ws.Cells.Item(22, 2).Value = 0.4374203038949512
$ws.Cells.Item(22, 3).Value = 0.2521171062027519
$ws.Cells.Item(22, 4).Value = 0.2974419321570869
$ws.Cells.Item(22, 5).Value = 2.542806401728791
$ws.Cells.Item(22, 6).Value = 10.0101393154344
$ws.Cells.Item(22, 7).Value = -8.161866677895206
$ws.Cells.Item(22, 8).Value = -0.02249986200433229
$ws.Cells.Item(22, 9).Value = 2.124945828975858
$ws.Cells.Item(22, 10).Value = -1.315135252606165
$ws.Cells.Item(22, 11).Value = -1.376330353346725
$ws.Cells.Item(23, 2).Value = 0.09526702093765275
$ws.Cells.Item(23, 3).Value = 0.2022742013914471
$ws.Cells.Item(23, 4).Value = 2.48837927106286
$ws.Cells.Item(23, 5).Value = 9.973585328022629
$ws.Cells.Item(23, 6).Value = -8.18953186637841
$ws.Cells.Item(23, 7).Value = -0.04590700167627143
$ws.Cells.Item(23, 8).Value = 2.10360890909983
$ws.Cells.Item(23, 9).Value = -1.335467850351198
$ws.Cells.Item(23, 10).Value = -1.396174332378672
$ws.Cells.Item(23, 11).Value = 0.6455109801188312
$ws.Cells.Item(24, 2).Value = 0.1534105877790901
$ws.Cells.Item(24, 3).Value = 2.460610420786189
$ws.Cells.Item(24, 4).Value = 9.956365482985138
$ws.Cells.Item(24, 5).Value = -8.201626982332003
$ws.Cells.Item(24, 6).Value = -0.05549112916037019
$ws.Cells.Item(24, 7).Value = 2.095251978282879
$ws.Cells.Item(24, 8).Value = -1.343224567749304
$ws.Cells.Item(24, 9).Value = -1.403637555454572
$ws.Cells.Item(24, 10).Value = 0.6381912792035592
$ws.Cells.Item(24, 11).Value = 0.05577943355819903
